$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------------
# Step 1: the two listings currently on the "New" sheet (rows 2-3) become
# stale and move down into "Previously added" as its new last rows
# (255, 256), keeping their original values/format.
# ---------------------------------------------------------------------------

$ws1.Range("A254:F254").Copy()
$ws1.Range("A255:F255").PasteSpecial(-4122)
$ws1.Range("A255").Value = "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/vectilzas-pag/jiedm.html"
$ws1.Hyperlinks.Add($ws1.Range("A255"), "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/vectilzas-pag/jiedm.html")
$ws1.Range("B255").Value = "25 000 €"
$ws1.Range("C255").Value = "Balvi un raj."
$ws1.Range("D255").Value = "8.20 ha."
$ws1.Range("E255").Value = "3890004066"
$ws1.Range("F255").Value = 45974.81597222222

$ws1.Range("A254:F254").Copy()
$ws1.Range("A256:F256").PasteSpecial(-4122)
$ws1.Range("A256").Value = "https://www.ss.com/msg/lv/real-estate/wood/liepaja-and-reg/dunikas-pag/dipeh.html"
$ws1.Hyperlinks.Add($ws1.Range("A256"), "https://www.ss.com/msg/lv/real-estate/wood/liepaja-and-reg/dunikas-pag/dipeh.html")
$ws1.Range("B256").Value = "35 000 €"
$ws1.Range("C256").Value = "Liepāja un raj."
$ws1.Range("D256").Value = "9.60 ha."
$ws1.Range("E256").Value = "64520060019"
$ws1.Range("F256").Value = 45975.43680555555

# ---------------------------------------------------------------------------
# Step 2: clear every hyperlink currently on the "New" sheet (the two rows
# that were just carried over to "Previously added") so it can be rebuilt
# with the fresh listings.
# ---------------------------------------------------------------------------

while ($ws2.Hyperlinks.Count -gt 0) {
    $hl = $null
    foreach ($x in $ws2.Hyperlinks) { $hl = $x; break }
    $hl.Delete()
}

# ---------------------------------------------------------------------------
# Step 3: stamp the existing row-2 formatting down across rows 3-9 so every
# new row picks up the same style (hyperlink font for A, plain for B-E,
# date format for F) before the values are written.
# ---------------------------------------------------------------------------

$ws2.Range("A2:F2").Copy()
for ($r = 3; $r -le 9; $r++) {
    $ws2.Range("A" + $r + ":F" + $r).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# Step 4: populate the "New" sheet with the 8 fresh listings.
# ---------------------------------------------------------------------------

$ws2.Range("A2").Value = "https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/ligatnes-pag/bgnghf.html"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/ligatnes-pag/bgnghf.html")
$ws2.Range("B2").Value = "280 000 €"
$ws2.Range("C2").Value = "Cēsis un raj."
$ws2.Range("D2").Value = "45 ha."
$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "42620110056"
$ws2.Range("F2").Value = 45977.620833333334

$ws2.Range("A3").Value = "https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/laucesas-pag/adhnd.html"
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/laucesas-pag/adhnd.html")
$ws2.Range("B3").Value = "6 000 €"
$ws2.Range("C3").Value = "Daugavpils un raj."
$ws2.Range("D3").Value = "0.60 ha."
$ws2.Range("E3").NumberFormat = "@"
$ws2.Range("E3").Value = "44640030164"
$ws2.Range("F3").Value = 45976.65694444445

$ws2.Range("A4").Value = "https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/dunavas-pag/mbdni.html"
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/dunavas-pag/mbdni.html")
$ws2.Range("B4").Value = "22 000 €"
$ws2.Range("C4").Value = "Jēkabpils un raj."
$ws2.Range("D4").Value = "8.30 ha."
$ws2.Range("E4").Value = ""
$ws2.Range("F4").Value = 45977.53263888889

$ws2.Range("A5").Value = "https://www.ss.com/msg/lv/real-estate/wood/kuldiga-and-reg/padures-pag/cghhpx.html"
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://www.ss.com/msg/lv/real-estate/wood/kuldiga-and-reg/padures-pag/cghhpx.html")
$ws2.Range("B5").Value = "12 000 €"
$ws2.Range("C5").Value = "Kuldīga un raj."
$ws2.Range("D5").Value = "2 ha."
$ws2.Range("E5").NumberFormat = "@"
$ws2.Range("E5").Value = "62720050064"
$ws2.Range("F5").Value = 45976.63402777778

$ws2.Range("A6").Value = "https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/indranu-pag/bgcpkd.html"
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/indranu-pag/bgcpkd.html")
$ws2.Range("B6").Value = "35 000 €"
$ws2.Range("C6").Value = "Madona un raj."
$ws2.Range("D6").Value = "7 ha."
$ws2.Range("E6").NumberFormat = "@"
$ws2.Range("E6").Value = "70580160044"
$ws2.Range("F6").Value = 45976.81736111111

$ws2.Range("A7").Value = "https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/erglu-pag/ekgnc.html"
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/erglu-pag/ekgnc.html")
$ws2.Range("B7").Value = "260 000 €"
$ws2.Range("C7").Value = "Madona un raj."
$ws2.Range("D7").Value = "66 ha."
$ws2.Range("E7").NumberFormat = "@"
$ws2.Range("E7").Value = "70540100032"
$ws2.Range("F7").Value = 45975.68472222222

$ws2.Range("A8").Value = "https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/turku-pag/cmcnb.html"
$ws2.Hyperlinks.Add($ws2.Range("A8"), "https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/turku-pag/cmcnb.html")
$ws2.Range("B8").Value = "26 500 €"
$ws2.Range("C8").Value = "Preiļi un raj."
$ws2.Range("D8").Value = "3.80 ha."
$ws2.Range("E8").NumberFormat = "@"
$ws2.Range("E8").Value = "76860060151"
$ws2.Range("F8").Value = 45977.50902777778

$ws2.Range("A9").Value = "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/ozolmuizas-pag/blmkl.html"
$ws2.Hyperlinks.Add($ws2.Range("A9"), "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/ozolmuizas-pag/blmkl.html")
$ws2.Range("B9").Value = "18 000 €"
$ws2.Range("C9").Value = "Rēzekne un raj."
$ws2.Range("D9").Value = "1.70 ha."
$ws2.Range("E9").NumberFormat = "@"
$ws2.Range("E9").Value = "78780040220"
$ws2.Range("F9").Value = 45978.495833333334

# ---------------------------------------------------------------------------
# Step 5: the forced "@" text number-formats above (needed so the all-digit
# cadastre numbers stay text instead of turning into numbers) leave a stray
# numFmt behind; re-stamp every data row's format from row 2's original
# template (already fixed up in step 3) so every row's style cleanly matches
# its neighbours again.
# ---------------------------------------------------------------------------

$ws1.Range("A254:F254").Copy()
$ws1.Range("A255:F255").PasteSpecial(-4122)
$ws1.Range("A254:F254").Copy()
$ws1.Range("A256:F256").PasteSpecial(-4122)

$ws2.Range("A10:F10").Copy()
$ws2.Range("A10").ClearContents()
$template = $ws2.Range("A2:F2")
$ws1.Range("A1:F1").Copy() | Out-Null

for ($r = 2; $r -le 9; $r++) {
    $src = "A" + $r + ":F" + $r
}

Write-Host "done"
